$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.907.43"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "1.551.47"
$ws.Range("E3").Value = "  -0.38%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("E7").Value = "  -0.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.20%  "

$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("E11").Value = "  -0.65%  "

$ws.Range("D12").Value = "1.773.20"
$ws.Range("E12").Value = "  -0.34%  "

$ws.Range("D13").Value = "1.547.74"
$ws.Range("E13").Value = "  -0.62%  "

$ws.Range("E14").Value = "  +0.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.518"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.44%  "

$ws.Range("D16").Value = "26.909.37"
$ws.Range("E16").Value = "  -0.36%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.64%  "

$ws.Range("E18").Value = "  +3.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("E20").Value = "  +0.11%  "

$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.18"
$ws.Range("D23").Style = "Normal"

$ws.Range("E24").Value = "  -1.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.53"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.69%  "

$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("E27").Value = "  +0.32%  "

$ws.Range("E28").Value = "  +0.39%  "

$ws.Range("E29").Value = "  -0.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0468"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.07%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("E33").Value = "  +3.57%  "

$ws.Range("D34").Value = "1.411.89"
$ws.Range("E34").Value = "  +0.75%  "

$ws.Range("E35").Value = "  +1.65%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.967"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("E37").Value = "  +0.20%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("E39").Value = "  +0.39%  "

$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("E41").Value = "  -0.54%  "

$ws.Range("E42").Value = "  +3.23%  "

$ws.Range("E43").Value = "  +1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.43"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("E46").Value = "  -1.12%  "

$ws.Range("D47").Value = "1.686.91"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0521"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.83%  "

$ws.Range("E50").Value = "  +2.90%  "

$ws.Range("E51").Value = "  -0.02%  "
